# The two "manager id / password" pairs used for rows 2 and 4 of the
# UserName sheet are refreshed to a new randomly-generated test value
# (this mirrors the "updated new test cases" commit).
#
# Old: A2/A4 = "mngr357958", B2/B4 = "rypabEz"
# New: A2/A4 = "mngr365881", B2/B4 = "jYmebUz"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "mngr365881"
$ws.Range("B2").Value = "jYmebUz"
$ws.Range("A4").Value = "mngr365881"
$ws.Range("B4").Value = "jYmebUz"

# The sheet's last active selection moved to J16 before the file was saved.
$ws.Range("J16").Select() | Out-Null
